$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Exponential distribution" results column is inserted before the
# existing "mm"/"sm"/"ms"/"ss" columns (which shift from C/D to D/E).
# Row 1 (the "service time distributions" banner in C1) is left untouched.

$ws.Range("C2").Value = "Exponential distribution"
$ws.Range("D2").Value = "mm"
$ws.Range("E2").Value = "sm"

# --- servers = 1 -----------------------------------------------------
$ws.Range("C3").Value = 52.947902323969799
$ws.Range("D3").Value = 225569.26222046599
$ws.Range("E3").Value = 149740.05733049399

$ws.Range("C4").Value = 44.162893299676398
$ws.Range("D4").Value = 129956.723060284
$ws.Range("E4").Value = 86382.683362601005

$ws.Range("C5").Value = 2.1911512597691098
$ws.Range("D5").Value = 538.79623599768797
$ws.Range("E5").Value = 358.26010415643702

$ws.Range("C6").Value = 0.93813258652195597
$ws.Range("D6").Value = 44.252924375031
$ws.Range("E6").Value = 26.193638724849801

# --- servers = 2 -----------------------------------------------------
$ws.Range("C7").Value = [double]"4.2279784780093102E-02"
$ws.Range("D7").Value = 104912.220439422
$ws.Range("E7").Value = 67722.966708097403

$ws.Range("C8").Value = [double]"1.0359324443449801E-02"
$ws.Range("D8").Value = 60642.989350169199
$ws.Range("E8").Value = 39004.312617370102

$ws.Range("C9").Value = [double]"9.0423864873150198E-02"
$ws.Range("D9").Value = 251.709669571646
$ws.Range("E9").Value = 162.44144147346401

$ws.Range("C10").Value = [double]"2.0946151852236999E-02"
$ws.Range("D10").Value = 22.553671252861701
$ws.Range("E10").Value = 12.9469997073386

# --- servers = 4 -----------------------------------------------------
$ws.Range("C11").Value = [double]"8.3081053291551095E-04"
$ws.Range("D11").Value = 44777.368940481203
$ws.Range("E11").Value = 26271.194566919199

$ws.Range("C12").Value = [double]"4.8583102563977399E-04"
$ws.Range("D12").Value = 25826.104313972599
$ws.Range("E12").Value = 15150.657471050999

$ws.Range("C13").Value = [double]"7.5328932305506599E-03"
$ws.Range("D13").Value = 107.704667387314
$ws.Range("E13").Value = 62.8980383635419

$ws.Range("C14").Value = [double]"3.3629020115662298E-03"
$ws.Range("D14").Value = 10.3373377318872
$ws.Range("E14").Value = 5.4421145594641098

# --- servers = 8 -------------------------------------------------------
# These exponential-distribution values are tiny, so (like the rest of the
# row) they get a scientific-notation number format.
$ws.Range("C15").Value = [double]"2.0469768189286601E-07"
$ws.Range("D15").Value = 14972.739586018501
$ws.Range("E15").Value = 5538.6466320352001

$ws.Range("C16").Value = [double]"2.8876161541123001E-06"
$ws.Range("D16").Value = 8690.4685143093502
$ws.Range("E16").Value = 3220.0367116898901

$ws.Range("C17").Value = [double]"6.4698717013045196E-06"
$ws.Range("D17").Value = 35.823409233631402
$ws.Range("E17").Value = 13.451919553566

$ws.Range("C18").Value = [double]"9.1268771912614605E-05"
$ws.Range("D18").Value = 4.8847667814621598
$ws.Range("E18").Value = 2.1979845201460302

$ws.Range("C15:C18").NumberFormat = "0.00E+00"

# New / adjusted column widths.
$ws.Columns("C:C").ColumnWidth = 23
$ws.Columns("F:F").ColumnWidth = 11

# The saved selection marker moved from D25 to C25.
$ws.Range("C25").Select()
